$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new notes next to the "Time spend" row data (column G),
# explaining how the numbers were derived.
$ws.Range("G2").Value = "<- Rounded sum of all values"
$ws.Range("G3").Value = "Calculated by hand. Survey tool of choice was not that advanced. My bad on that one. "

# Nudge "Chart 1" (the Login/Register chart) a bit to the left and down,
# matching how it was manually repositioned.
$co = $ws.ChartObjects(1)
$co.Left = $co.Left - 18.75
$co.Top = $co.Top + 4.5

# Update the saved selection to reflect where the user left the cursor.
$ws.Range("S5").Select() | Out-Null
